# Update the "想去人数" (want-to-go count) figures in column F of the
# "展览" and "全部类型" sheets for rows 2-6.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value mapping for column F
$updates = @{
    2 = 6348
    3 = 26
    4 = 186
    5 = 1009
    6 = 107
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
